$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 619.6923
$ws.Range("I33").Value = 663
$ws.Range("J33").Value = 100
$ws.Range("K33").Value = 663
$ws.Range("L33").Value = 100
$ws.Range("M33").Value = -434
$ws.Range("N33").Value = -558
# Row 55
$ws.Range("H55").Value = 71429630
$ws.Range("I55").Value = 90910320
$ws.Range("J55").Value = 467.33334
$ws.Range("K55").Value = 90910320
$ws.Range("L55").Value = 467.33334
$ws.Range("M55").Value = -90910106
$ws.Range("N55").Value = -895.33334
# Row 94
$ws.Range("H94").Value = 2128.3333
$ws.Range("I94").Value = 2128.3333
$ws.Range("K94").Value = 2128.3333
$ws.Range("M94").Value = -1677.3333
# Row 132
$ws.Range("H132").Value = 3201.84
$ws.Range("I132").Value = 3103.35
$ws.Range("J132").Value = 3595.8
$ws.Range("K132").Value = 9310.049999999999
$ws.Range("L132").Value = 10787.4
$ws.Range("M132").Value = -6780.049999999999
$ws.Range("N132").Value = -15847.4
# Row 137
$ws.Range("H137").Value = 1257.5111
$ws.Range("I137").Value = 936
$ws.Range("J137").Value = 1840.25
$ws.Range("K137").Value = 2808
$ws.Range("L137").Value = 5520.75
$ws.Range("M137").Value = -258
$ws.Range("N137").Value = -10620.75

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 984.087
$ws.Range("I2").Value = 935.2222
$ws.Range("J2").Value = 1160
$ws.Range("K2").Value = 935.2222
$ws.Range("L2").Value = 1160
$ws.Range("M2").Value = -822.2222
$ws.Range("N2").Value = -1386
# Row 5
$ws.Range("H5").Value = 133.85715
$ws.Range("I5").Value = 67.40000000000001
$ws.Range("J5").Value = 300
$ws.Range("K5").Value = 67.40000000000001
$ws.Range("L5").Value = 300
$ws.Range("M5").Value = 44.59999999999999
$ws.Range("N5").Value = -524
# Row 45
$ws.Range("H45").Value = 2066.0557
$ws.Range("I45").Value = 1838.4615
$ws.Range("J45").Value = 2657.8
$ws.Range("K45").Value = 1838.4615
$ws.Range("L45").Value = 2657.8
$ws.Range("M45").Value = -1461.4615
$ws.Range("N45").Value = -3411.8
# Row 74
$ws.Range("H74").Value = 2171.6292
$ws.Range("I74").Value = 1224.4324
$ws.Range("J74").Value = 3573.48
$ws.Range("K74").Value = 1224.4324
$ws.Range("L74").Value = 3573.48
$ws.Range("M74").Value = -350.4323999999999
$ws.Range("N74").Value = -5321.48
# Row 77
$ws.Range("H77").Value = 2171.6292
$ws.Range("I77").Value = 1224.4324
$ws.Range("J77").Value = 3573.48
$ws.Range("K77").Value = 6122.161999999999
$ws.Range("L77").Value = 17867.4
$ws.Range("M77").Value = -1754.161999999999
$ws.Range("N77").Value = -26603.4
# Row 116
$ws.Range("H116").Value = 984.087
$ws.Range("I116").Value = 935.2222
$ws.Range("J116").Value = 1160
$ws.Range("K116").Value = 935.2222
$ws.Range("L116").Value = 1160
$ws.Range("M116").Value = 1358.7778
$ws.Range("N116").Value = -5748
# Row 132
$ws.Range("H132").Value = 5979.353
$ws.Range("I132").Value = 4465.2
$ws.Range("J132").Value = 8142.4287
$ws.Range("K132").Value = 13395.6
$ws.Range("L132").Value = 24427.2861
$ws.Range("M132").Value = -10865.6
$ws.Range("N132").Value = -29487.2861

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 984.087
$ws.Range("I3").Value = 935.2222
$ws.Range("J3").Value = 1160
$ws.Range("K3").Value = 935.2222
$ws.Range("L3").Value = 1160
$ws.Range("M3").Value = -821.2222
$ws.Range("N3").Value = -1388
# Row 4
$ws.Range("H4").Value = 133.85715
$ws.Range("I4").Value = 67.40000000000001
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 67.40000000000001
$ws.Range("L4").Value = 300
$ws.Range("M4").Value = 47.59999999999999
$ws.Range("N4").Value = -530
# Row 22
$ws.Range("H22").Value = 416.66666
$ws.Range("I22").Value = 416.66666
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 416.66666
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -243.66666
$ws.Range("N22").ClearContents()
# Row 99
$ws.Range("H99").Value = 2303.875
$ws.Range("I99").Value = 2006.4
$ws.Range("K99").Value = 2006.4
$ws.Range("M99").Value = -508.4000000000001

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 62.142857
$ws.Range("I7").Value = 37
$ws.Range("K7").Value = 37
$ws.Range("M7").Value = 76
# Row 58
$ws.Range("H58").Value = 1768.5254
$ws.Range("I58").Value = 1233.5476
$ws.Range("J58").Value = 3090.2354
$ws.Range("K58").Value = 1233.5476
$ws.Range("L58").Value = 3090.2354
$ws.Range("M58").Value = -1030.5476
$ws.Range("N58").Value = -3496.2354
# Row 132
$ws.Range("H132").Value = 6669027.5
$ws.Range("I132").Value = 2265.9285
$ws.Range("J132").Value = 15153997
$ws.Range("K132").Value = 6797.7855
$ws.Range("L132").Value = 45461991
$ws.Range("M132").Value = -4267.7855
$ws.Range("N132").Value = -45467051
# Row 134
$ws.Range("H134").Value = 7619.8096
$ws.Range("I134").Value = 7211.1055
$ws.Range("J134").Value = 11502.5
$ws.Range("K134").Value = 21633.3165
$ws.Range("L134").Value = 34507.5
$ws.Range("M134").Value = -19098.3165
$ws.Range("N134").Value = -39577.5
# Row 136
$ws.Range("H136").Value = 1768.5254
$ws.Range("I136").Value = 1233.5476
$ws.Range("J136").Value = 3090.2354
$ws.Range("K136").Value = 3700.642800000001
$ws.Range("L136").Value = 9270.706200000001
$ws.Range("M136").Value = -1150.642800000001
$ws.Range("N136").Value = -14370.7062

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 637.871
$ws.Range("I5").Value = 489.07144
$ws.Range("J5").Value = 2026.6666
$ws.Range("K5").Value = 1467.21432
$ws.Range("L5").Value = 6079.9998
$ws.Range("M5").Value = -1355.21432
$ws.Range("N5").Value = -6303.9998
# Row 31
$ws.Range("H31").Value = 2424.9375
$ws.Range("I31").Value = 1599.5
$ws.Range("J31").Value = 2542.8572
$ws.Range("K31").Value = 4798.5
$ws.Range("L31").Value = 7628.571599999999
$ws.Range("M31").Value = -4510.5
$ws.Range("N31").Value = -8204.571599999999
# Row 75
$ws.Range("H75").Value = 1599.6666
$ws.Range("J75").Value = 2999
$ws.Range("L75").Value = 8997
$ws.Range("N75").Value = -10993
# Row 78
$ws.Range("H78").Value = 1599.6666
$ws.Range("J78").Value = 2999
$ws.Range("L78").Value = 26991
$ws.Range("N78").Value = -36975
# Row 102
$ws.Range("H102").Value = 3500
$ws.Range("J102").Value = 3500
$ws.Range("L102").Value = 10500
$ws.Range("N102").Value = -15368
# Row 104
$ws.Range("H104").Value = 5243
$ws.Range("J104").Value = 5243
$ws.Range("L104").Value = 15729
$ws.Range("N104").Value = -20971
# Row 131
$ws.Range("H131").Value = 4013.4285
$ws.Range("J131").Value = 4839.8237
$ws.Range("L131").Value = 14519.4711
$ws.Range("N131").Value = -24599.4711
# Row 135
$ws.Range("H135").Value = 637.871
$ws.Range("I135").Value = 489.07144
$ws.Range("J135").Value = 2026.6666
$ws.Range("K135").Value = 4401.64296
$ws.Range("L135").Value = 18239.9994
$ws.Range("M135").Value = -1866.64296
$ws.Range("N135").Value = -23309.9994

$ws = $wb.Worksheets.Item("GSM")
# Row 113
$ws.Range("H113").Value = 47325.582
$ws.Range("I113").Value = 59321.79
$ws.Range("J113").Value = 1740
$ws.Range("K113").Value = 59321.79
$ws.Range("L113").Value = 1740
$ws.Range("M113").Value = -57151.79
$ws.Range("N113").Value = -6080
# Row 126
$ws.Range("H126").Value = 2217.5293
$ws.Range("I126").Value = 2184.5715
$ws.Range("K126").Value = 6553.7145
$ws.Range("M126").Value = -4083.7145
# Row 132
$ws.Range("H132").Value = 3161.6
$ws.Range("I132").Value = 2505.36
$ws.Range("J132").Value = 4802.2
$ws.Range("K132").Value = 7516.08
$ws.Range("L132").Value = 14406.6
$ws.Range("M132").Value = -4986.08
$ws.Range("N132").Value = -19466.6

$ws = $wb.Worksheets.Item("LTW")
# Row 136
$ws.Range("H136").Value = 11112886
$ws.Range("I136").Value = 1610.3636
$ws.Range("J136").Value = 41668892
$ws.Range("K136").Value = 4831.0908
$ws.Range("L136").Value = 125006676
$ws.Range("M136").Value = -2281.0908
$ws.Range("N136").Value = -125011776
